$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Notes added in column D describing a power-outage/requeue event that affected
# several queued runs. Each "run" block (date / last restart avail / gmax /
# hotstart / coldstart / initC) gets short free-text notes next to most rows.

$ws.Range("D3").Value = "22 june"
$ws.Range("D4").Value = "power "
$ws.Range("D5").Value = "out"
$ws.Range("D6").Value = "requeued"
$ws.Range("D7").Value = "auto"

$ws.Range("D11").Value = "22 june"
$ws.Range("D12").Value = "power "
$ws.Range("D13").Value = "out"
$ws.Range("D14").Value = "requeued"
$ws.Range("D15").Value = "auto"

$ws.Range("D19").Value = "22 june"
$ws.Range("D20").Value = "power "
$ws.Range("D21").Value = "out"
$ws.Range("D22").Value = "requeued"
$ws.Range("D23").Value = "auto"

$ws.Range("D27").Value = "22 june"
$ws.Range("D28").Value = "power "
$ws.Range("D29").Value = "out"
$ws.Range("D30").Value = "requeued"
$ws.Range("D31").Value = "auto"

$ws.Range("D35").Value = "22 june"
$ws.Range("D36").Value = "power out"
$ws.Range("D37").Value = "requeued"
$ws.Range("D38").Value = "auto"

$ws.Range("D43").Value = "22 june"
$ws.Range("D44").Value = "power"
$ws.Range("D45").Value = "out"
$ws.Range("D46").Value = "requeued "
$ws.Range("D47").Value = "auto"

$ws.Range("D51").Value = "22 june"
$ws.Range("D52").Value = 1516

# "0.50" looks numeric, so Excel would normally coerce it to the number 0.5.
# Force it to be stored as text (matching the original author's cell), then
# restore the plain "Normal" style so no extra number-format style lingers.
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "0.50"
$ws.Range("D53").Style = "Normal"

$ws.Range("D54").Value = "no"
$ws.Range("D55").Value = "yes"
$ws.Range("D56").Value = 1516

# Match the final selection left by the author: D3:D7 active cell D3.
$ws.Range("D3:D7").Select()
